# Generate Report for handoff
#
# The previous handoff attempt (a237ea19-6606-4ed2-9d64-215b69acb2df.md) is
# superseded by a new source file (2969e9b0-6ac0-4953-b192-c79b370797b6.md)
# which successfully handed off, while a second source file
# (e860b43e-532b-4a59-8294-a8b42b6c6967.md) failed transform, and the
# .localization-config file (previously "not to be localized" in row 3) is
# now tracked on its own new row.

$wb = $excel.ActiveWorkbook

# ---- shared constants -------------------------------------------------
$mdNew        = "2969e9b0-6ac0-4953-b192-c79b370797b6.md"
$mdFailed     = "e860b43e-532b-4a59-8294-a8b42b6c6967.md"
$configName   = ".localization-config"

$readyStatus   = "Ready for handoff"
$failedStatus  = "Handoff transform failed"
$notLocStatus  = "Not to be localized"

$baseCommit   = "fa2b12bc368f421dfa17dabf84757281f70e3b63"
$mdNewUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/$baseCommit/e2e/$mdNew"
$mdFailedUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/$baseCommit/e2e/$mdFailed"
$configUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/$baseCommit/$configName"

$mdNewStem    = "2969e9b0-6ac0-4953-b192-c79b370797b6"
$zhHash       = "d875920e262500ef34361872380c4b4e53072503"
$deHash       = "d875920e262500ef34361872380c4b4e53072503"
$zhXlf        = "$mdNewStem.$zhHash.zh-cn.xlf"
$deXlf        = "$mdNewStem.$deHash.de-de.xlf"

$zhCommit     = "a664c960643954fc79cdd24885679a4b6f78f57b"
$deCommit     = "a72b6fd5080ab13d3e476413d8b3fee3d3ca487e"
$zhXlfUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/$zhXlf"
$deXlfUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/$deXlf"

$zhHandoffDt  = "2016-01-18 03:08:32"
$deHandoffDt  = "2016-01-18 03:08:43"
$epochDt      = "0001-01-01 00:00:00"
$dtFormat     = "yyyy-mm-dd HH:mm:ss"

# =========================================================================
# Sheet 1: Overview
# =========================================================================
$ws1 = $wb.Worksheets.Item("Overview")

# Drop the stale hyperlink registrations so we can cleanly re-add them
# against the new targets (row 1 has no hyperlink, so this is a safe anchor).
$ws1.Range("A1").Hyperlinks.Delete()

$ws1.Hyperlinks.Add($ws1.Range("A2"), $mdNewUrl, "", "", $mdNew)
$ws1.Range("B2").Value = $readyStatus
$ws1.Range("C2").Value = $readyStatus

$ws1.Hyperlinks.Add($ws1.Range("A3"), $mdFailedUrl, "", "", $mdFailed)
$ws1.Range("B3").Value = $failedStatus
$ws1.Range("C3").Value = $failedStatus

$ws1.Hyperlinks.Add($ws1.Range("A4"), $configUrl, "", "", $configName)
$ws1.Range("B4").Value = $notLocStatus
$ws1.Range("C4").Value = $notLocStatus

# =========================================================================
# Sheet 2: zh-cn
# =========================================================================
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A1").Hyperlinks.Delete()

$ws2.Hyperlinks.Add($ws2.Range("A2"), $mdNewUrl, "", "", $mdNew)
$ws2.Range("B2").Value = $readyStatus
$ws2.Hyperlinks.Add($ws2.Range("C2"), $zhXlfUrl, "", "", $zhXlf)
$ws2.Range("D2").Value = $zhHandoffDt
$ws2.Range("G2").Value = $epochDt
$ws2.Range("H2").Value = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A3"), $mdFailedUrl, "", "", $mdFailed)
$ws2.Range("B3").Value = $failedStatus
$ws2.Range("D3").Value = $epochDt
$ws2.Range("G3").Value = $epochDt
$ws2.Range("H3").Value = "Ignored"

$ws2.Hyperlinks.Add($ws2.Range("A4"), $configUrl, "", "", $configName)
$ws2.Range("B4").Value = $notLocStatus
$ws2.Range("D4").Value = $epochDt
$ws2.Range("D4").NumberFormat = $dtFormat
$ws2.Range("G4").Value = $epochDt
$ws2.Range("H4").Value = "Ignored"

# =========================================================================
# Sheet 3: de-de
# =========================================================================
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A1").Hyperlinks.Delete()

$ws3.Hyperlinks.Add($ws3.Range("A2"), $mdNewUrl, "", "", $mdNew)
$ws3.Range("B2").Value = $readyStatus
$ws3.Hyperlinks.Add($ws3.Range("C2"), $deXlfUrl, "", "", $deXlf)
$ws3.Range("D2").Value = $deHandoffDt
$ws3.Range("G2").Value = $epochDt
$ws3.Range("H2").Value = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A3"), $mdFailedUrl, "", "", $mdFailed)
$ws3.Range("B3").Value = $failedStatus
$ws3.Range("D3").Value = $epochDt
$ws3.Range("G3").Value = $epochDt
$ws3.Range("H3").Value = "Ignored"

$ws3.Hyperlinks.Add($ws3.Range("A4"), $configUrl, "", "", $configName)
$ws3.Range("B4").Value = $notLocStatus
$ws3.Range("D4").Value = $epochDt
$ws3.Range("D4").NumberFormat = $dtFormat
$ws3.Range("G4").Value = $epochDt
$ws3.Range("H4").Value = "Ignored"

Write-Output "Report regenerated for handoff"
